$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The runs/balls/fours/sixes columns (C:F) store numbers as text in this
# sheet. Force text format on each changed row before assigning so the
# values keep their text type, matching the original workbook.

# Row 2
$ws.Range("C2:F2").NumberFormat = "@"
$ws.Range("C2").Value = "22"
$ws.Range("D2").Value = "12"
$ws.Range("E2").Value = "3"
$ws.Range("F2").Value = "1"

# Row 4
$ws.Range("C4:F4").NumberFormat = "@"
$ws.Range("C4").Value = "4"
$ws.Range("D4").Value = "7"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "0"

# Row 5
$ws.Range("C5:F5").NumberFormat = "@"
$ws.Range("C5").Value = "16"
$ws.Range("D5").Value = "13"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "1"

# Row 6
$ws.Range("C6:F6").NumberFormat = "@"
$ws.Range("C6").Value = "13"
$ws.Range("D6").Value = "8"
$ws.Range("E6").Value = "2"
$ws.Range("F6").Value = "0"

# Row 7
$ws.Range("C7:F7").NumberFormat = "@"
$ws.Range("C7").Value = "70"
$ws.Range("D7").Value = "44"
$ws.Range("E7").Value = "4"
$ws.Range("F7").Value = "5"

# Row 8
$ws.Range("C8:F8").NumberFormat = "@"
$ws.Range("C8").Value = "21"
$ws.Range("D8").Value = "16"
$ws.Range("E8").Value = "1"
$ws.Range("F8").Value = "2"
